# Applies the "Added -solution to project description doc" edit:
#  1. Removes the stray _GoBack bookmark that sits at the end of the
#     "Problem" heading paragraph.
#  2. Appends a new "Solution" heading paragraph (bold, matching the
#     "Problem" heading formatting) plus a body paragraph with the
#     solution text, right after the existing problem-description
#     paragraph. The _GoBack bookmark is re-created in the middle of
#     the new body paragraph (Word always parks it at the last edit
#     point), splitting the text at "...hence providi|ng a secure
#     channel.".

$d = $word.ActiveDocument

# --- 1. Drop the old _GoBack bookmark --------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Create the two new (still empty) paragraphs first ------------------
# Doing this before any bold formatting is applied means the paragraph mark
# each new paragraph inherits stays non-bold (matching the "Wireless
# devices ..." paragraph it grew out of), so the body paragraph never picks
# up an explicit Bold flag later on.
$bodyPara = $d.Paragraphs.Item(2)            # "Wireless devices ..." paragraph
$insertPoint = $bodyPara.Range
$insertPoint.Collapse(0)                     # 0 = wdCollapseEnd
$insertPoint.InsertParagraphAfter()          # -> paragraph 3 (heading, empty)

$headingEnd = $d.Paragraphs.Item(3).Range
$headingEnd.Collapse(0)
$headingEnd.InsertParagraphAfter()           # -> paragraph 4 (body, empty)

# --- 3. Fill in + format the "Solution" heading paragraph ------------------
$headingRange = $d.Paragraphs.Item(3).Range
$headingRange.InsertAfter("Solution")
$headingRange = $d.Paragraphs.Item(3).Range
$headingRange.Font.Name = "Times New Roman"
$headingRange.Font.Bold = $true
$headingRange.Font.Size = 12

# --- 4. Fill in the solution body paragraph ---------------------------------
$solutionText = "To solve this problem, this proof-of-concept project proposes a system which uses a low frequency signal alongside magnetic fields passing through the body (as a channel) to enable communication between any two devices installed on the human body. Some of the most important advantages of this system is that, the information relayed in form of magnetic fields is capable of penetrating freely through the body tissues and therefore the communication can be achieved with lower signal loss (as a result of obstruction), and less power consumption. Besides, by using the human body as a channel of communication, the medium is less vulnerable to leak information in that, the data transmitted from one part of the body to another cannot be radiated off the body hence providing a secure channel."

$solutionRange = $d.Paragraphs.Item(4).Range
$solutionRange.InsertAfter($solutionText)
$solutionRange = $d.Paragraphs.Item(4).Range
$solutionRange.Font.Name = "Times New Roman"
$solutionRange.Font.Size = 12

# --- 5. Re-plant the _GoBack bookmark mid-paragraph (splits "providi"/"ng") -
$solutionPara = $d.Paragraphs.Item(4).Range
$splitPoint = $solutionPara.Start + $solutionText.IndexOf("hence providi") + "hence providi".Length
$goBackRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
